$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for first row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 01:02:17"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for first row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-19 01:02:12"
$wsZhCn.Range("K2").Value = "2016-08-19 01:02:30"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for first row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-19 01:02:17"
$wsDeDe.Range("K2").Value = "2016-08-19 01:02:37"
